# Apply the "Fuel Price" column addition to the Engine Data sheet, fix a
# typo'd Mass value, and restore the saved cell selections.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Engine Data")
$ws2 = $wb.Worksheets.Item("Propellant Data")

# New column header
$ws1.Cells.Item(1, 10).Value = "Fuel Price"

# Fuel price per row, keyed by the propellant used in that row:
#   L02/ Kerosene -> 75.12  (matches the "75.12 p/L" note on Propellant Data)
#   L02/RP1       -> 93.87  (matches the "$93.87 per gallon" note)
#   anything else -> 1
$fuelPrices = @{
    2  = 75.12
    3  = 75.12
    4  = 75.12
    5  = 93.87
    6  = 93.87
    7  = 93.87
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 75.12
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
}

foreach ($row in $fuelPrices.Keys) {
    $ws1.Cells.Item($row, 10).Value = $fuelPrices[$row]
}

# Fix a typo in the Mass column for row 15 (RS-68): 6.597 -> 6597
$ws1.Cells.Item(15, 9).Value = 6597

# Restore the cell selections recorded in the saved workbook
$ws1.Activate()
[void]$ws1.Range("I15").Select()

$ws2.Activate()
[void]$ws2.Range("B4").Select()

$ws1.Activate()
